# Update cryptocurrency price/volume data as scraped on Fri Dec 22 21:41:48 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.833.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.317.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.85%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "97.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "272.65"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.42"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.04"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.31%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.659.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.15%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.51%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +8.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.326.14"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.773.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.37"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.29%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.29%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.50"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.43"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0918"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.49"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.62%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0360"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.78%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.245"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.18%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +21.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.21"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.85"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +9.97%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.50"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.195"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +17.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.545.81"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.44%  "
